$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (Algerian Ligue 1, MC Alger vs ES Ben Aknoun)
$ws.Range("F3").Value  = 1.34
$ws.Range("G3").Value  = 1.63
$ws.Range("H3").Value  = 8
$ws.Range("J3").Value  = 2.6
$ws.Range("K3").Value  = 7.8
$ws.Range("L3").Value  = 1.01
$ws.Range("M3").Value  = 1.01
$ws.Range("N3").Value  = 1.5
$ws.Range("O3").Value  = 1.01
$ws.Range("P3").Value  = 1.5
$ws.Range("Q3").Value  = 2.06
$ws.Range("R3").Value  = 1.12
$ws.Range("S3").Value  = 2.06
$ws.Range("T3").Value  = 1.01
$ws.Range("U3").Value  = 1.01
$ws.Range("V3").Value  = 1.01
$ws.Range("W3").Value  = 2.58
$ws.Range("X3").Value  = 1000
$ws.Range("Y3").Value  = 1000
$ws.Range("Z3").Value  = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 4 (Portuguese Primeira Liga, Guimaraes vs Sporting Lisbon)
$ws.Range("H4").Value  = 1.39
$ws.Range("K4").Value  = 5.4
$ws.Range("N4").Value  = 4
$ws.Range("O4").Value  = 1.29
$ws.Range("P4").Value  = 2.06
$ws.Range("Q4").Value  = 1.87
$ws.Range("R4").Value  = 1.41
$ws.Range("S4").Value  = 3.2
$ws.Range("T4").Value  = 2.08
$ws.Range("U4").Value  = 1.79
$ws.Range("V4").Value  = 3.25
$ws.Range("Y4").Value  = 10
$ws.Range("Z4").Value  = 7.8
$ws.Range("AB4").Value = 34
$ws.Range("AC4").Value = 13
$ws.Range("AD4").Value = 10
$ws.Range("AH4").Value = 32
$ws.Range("AJ4").Value = 470
$ws.Range("AK4").Value = 210
$ws.Range("AL4").Value = 180
$ws.Range("AN4").Value = 340
